$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-4.548879146575928, 1.983122944831848, -2.030394077301025),
    @(0.2330931723117828, 7.29985761642456, -8.993494987487793),
    @(4.196409702301025, -7.729958534240723, -5.114867210388184),
    @(1.506084203720093, -3.158908605575561, -2.592853784561157),
    @(0.7066032290458679, -5.12858247756958, 0.3414172530174255),
    @(-3.361641645431519, 6.300107002258301, 2.977946758270264),
    @(-2.117412805557251, 13.40914916992188, -1.111603498458862),
    @(1.100218415260315, 1.102149248123169, -2.866626739501953),
    @(0.4749082028865814, -12.12523937225342, -7.235276222229004),
    @(3.548196077346802, 9.086838722229004, 1.108407616615295)
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
